$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# --- Update row 28: fill previously-empty cells with the literal "nan" placeholder ---
$ws.Cells.Item(28, 1).Value = "nan"   # A28
$ws.Cells.Item(28, 12).Value = "nan"  # L28
$ws.Cells.Item(28, 13).Value = "nan"  # M28
$ws.Cells.Item(28, 15).Value = "nan"  # O28

# --- Add new row 29 describing the new event added to Card18 ---
# A29 must stay a text value ("18"), not be auto-converted to a number
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "18"                     # A29
$ws.Cells.Item(29, 1).ClearFormats()

$ws.Cells.Item(29, 12).Value = "15/1/2026"             # L29
$ws.Cells.Item(29, 13).Value = "قطع سير700"            # M29
$ws.Cells.Item(29, 15).Value = "تم تغير سير700"        # O29
$ws.Cells.Item(29, 16).Value = "مصطفي"                 # P29
